$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted price strings (e.g. "21.662.77"); force text
# format before assignment so Excel does not auto-convert them to numbers
# and strip formatting (trailing zeros, thousands dots, etc.).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "21.662.77"
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.532.69"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.00"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3929"
$ws.Range("E7").Value = "  +3.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3156"
$ws.Range("E8").Value = "  -2.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.36"
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07154"
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.043"
$ws.Range("E11").Value = "  -6.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.614"
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.45"
$ws.Range("E14").Value = "  -4.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.588"
$ws.Range("E15").Value = "  -3.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.538.25"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001091"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06589"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.82"
$ws.Range("E19").Value = "  -2.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.086"
$ws.Range("E21").Value = "  -5.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.37"
$ws.Range("E22").Value = "  -3.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.80"
$ws.Range("E23").Value = "  -5.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.383"
$ws.Range("E24").Value = "  +3.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "21.674.61"
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.324"
$ws.Range("E26").Value = "  -7.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.41"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.24"
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.833"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.712.81"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "116.99"
$ws.Range("E31").Value = "  -2.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.842"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9509"
$ws.Range("E33").Value = "  -13.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08107"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.571"
$ws.Range("E35").Value = "  -7.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06048"
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.079"
$ws.Range("E37").Value = "  -3.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02197"
$ws.Range("E38").Value = "  -3.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.448"
$ws.Range("E39").Value = "  -11.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2015"
$ws.Range("E40").Value = "  -4.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.177"
$ws.Range("E41").Value = "  -3.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5720"
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.02"
$ws.Range("E45").Value = "  -3.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.718"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5471"
$ws.Range("E47").Value = "  -4.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.152"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.857"
$ws.Range("E49").Value = "  -3.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "115.16"
$ws.Range("E50").Value = "  -3.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06688"
$ws.Range("E51").Value = "  -2.62%  "

# Rows 42/43: Aptos and Frax swapped positions (rank change) with updated data
$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.82"
$ws.Range("E43").Value = "  -0.22%  "
